# Request resource. Add conditional messaging. Style client rendering of preview controls
#
# The CAVEATS paragraph currently reads:
#   "To be revisited during the process of API utility."
# It must become two runs (same run formatting) reading:
#   "Add condition for when the preview URL " + "endpoint does not exist. Indicate to
#   end user that the resource is not available."

$d = $word.ActiveDocument

$part1 = "Add condition for when the preview URL "
$part2 = "endpoint does not exist. Indicate to end user that the resource is not available."

$rng = $d.Content
$found = $rng.Find.Execute("To be revisited during the process of API utility.", $true, $false, $false, $false, $false, $true, 1, $false, ($part1 + $part2), 2)

if ($found) {
    # $rng now spans the freshly inserted replacement text. Work out where the
    # boundary between the two desired runs should fall.
    $splitAt = $rng.Start + $part1.Length

    # Forcing a genuine run break (rather than one coalesced run) requires the
    # two halves to differ at the XML level for a moment. Dropping a bookmark at
    # the split point, then removing it again, leaves the surrounding run
    # formatting completely untouched but still splits the text into two
    # <w:r> runs with identical rPr - matching how the source document stores
    # the two sentences as separate runs.
    $markerName = "tmpRunSplitMarker"
    $markerRange = $d.Range($splitAt, $splitAt)
    $d.Bookmarks.Add($markerName, $markerRange)
    $d.Bookmarks($markerName).Delete()
}
